$wb = $excel.ActiveWorkbook
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Sheet "VENTAS POR GRUPO": zero-out several category cells, and refresh the
# "x de 12" counters in row 14 accordingly ---
$wsGrupo.Range("H3").Value = 0
$wsGrupo.Range("M9").Value = 0
$wsGrupo.Range("M10").Value = 0
$wsGrupo.Range("M11").Value = 0
$wsGrupo.Range("P11").Value = 0
$wsGrupo.Range("H12").Value = 0

$wsGrupo.Range("H14").Value = "0 de 12"
$wsGrupo.Range("M14").Value = "0 de 12"
$wsGrupo.Range("P14").Value = "0 de 12"

# --- Sheet "VENTA MENSUAL": monthly rollover. Column headers (junio..septiembre)
# advance by one month, the data in C..F shifts left by one column, and the new
# rightmost month (F) starts at zero. Column widths rotate the same way. ---

$wsMensual.Range("C1").Value = "julio"
$wsMensual.Range("D1").Value = "agosto"
$wsMensual.Range("E1").Value = "septiembre"
$wsMensual.Range("F1").Value = "octubre"

$oldColWidths = @{}
foreach ($col in @("C","D","E","F")) {
    $oldColWidths[$col] = $wsMensual.Range("$col`1").EntireColumn.ColumnWidth
}
$wsMensual.Range("C1").EntireColumn.ColumnWidth = $oldColWidths["D"]
$wsMensual.Range("D1").EntireColumn.ColumnWidth = $oldColWidths["E"]
$wsMensual.Range("E1").EntireColumn.ColumnWidth = $oldColWidths["F"]
$wsMensual.Range("F1").EntireColumn.ColumnWidth = $oldColWidths["C"]

$rows = @(3, 9, 10, 11, 12, 14)
foreach ($r in $rows) {
    $oldC = $wsMensual.Range("C$r").Value2
    $oldD = $wsMensual.Range("D$r").Value2
    $oldE = $wsMensual.Range("E$r").Value2
    $oldF = $wsMensual.Range("F$r").Value2

    $wsMensual.Range("C$r").Value = $oldD
    $wsMensual.Range("D$r").Value = $oldE
    $wsMensual.Range("E$r").Value = $oldF
    $wsMensual.Range("F$r").Value = 0
}
